$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old opportunity-tracker columns (B:L), keep column A for the
# new "Voucher Code" list. The previously-styled (but now empty) cells in
# E/F/J survive this as format-only cells, matching the target sheet.
$ws.Range("B1:L3").ClearContents()

# Replace the header/body text with the new voucher upload content.
$ws.Range("A1").Value = "Voucher Code"
$ws.Range("A2").Value = "Gyxwyzx123"
$ws.Range("A3").Value = "Vwyx55rXLt"

# The old WebLink hyperlinks (J2/J3) no longer apply - drop them.
$ws.Hyperlinks.Delete()

# Move the selection / top-left cell to A3 (also clears the stale
# topLeftCell="B1" scroll position left over from the old sheet).
$ws.Range("A3").Select()
